# Append: 2025-11-13 18:33 JST
# Update the "取得日時" (acquired timestamp) column A for all existing data
# rows on the "ランサーズ" sheet from the old run timestamp to the new one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$oldTimestamp = "2025-11-13 18:25:58"
$newTimestamp = "2025-11-13 18:33:41"

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) {
    $lastRow = 13
}

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    if ($cell.Text -eq $oldTimestamp) {
        $cell.Value = $newTimestamp
    }
}
